$d = $word.ActiveDocument

function Replace-AllText($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-AllText "0.9812 and coefficients" "1.0307 and coefficients"
Replace-AllText "2.9972, and" "3.0293, and"
Replace-AllText "1.9663" "2.02"

Replace-AllText "0.9808 and coefficients" "1.0298 and coefficients"
Replace-AllText "2.9966, and" "3.0289, and"
Replace-AllText "1.966" "2.0193"

Replace-AllText "array([0.95452431]) and coefficients" "array([1.01741031]) and coefficients"
Replace-AllText "array([3.01674673]), and" "array([3.03917656]), and"
Replace-AllText "array([1.99921879])" "array([2.03102484])"
